$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 for "RF" (random forest), pushing the existing "Ensemble" row to row 7.
# Copy formatting from row 5 so the new row matches the label-column style (bold, bordered, centered).
$ws.Rows.Item(6).Insert()
$ws.Range("A5:M5").Copy()
$ws.Range("A6:M6").PasteSpecial(-4122)

# Row labels
$ws.Range("A6").Value = "RF"

# Updated timing values for rows 2-5 (kNN, SVM, LR, NB) plus new row 6 (RF) and row 7 (Ensemble)
$ws.Range("B2").Value = 0.01223330497741699
$ws.Range("C2").Value = 0.02341761589050293
$ws.Range("D2").Value = 0.00535745620727539
$ws.Range("E2").Value = 0.01214289665222168
$ws.Range("F2").Value = 0.003723478317260742
$ws.Range("G2").Value = 0.05266880989074707
$ws.Range("H2").Value = 0.01597180366516113
$ws.Range("I2").Value = 0.01503376960754394
$ws.Range("J2").Value = 0.01118769645690918
$ws.Range("K2").Value = 0.01879720687866211
$ws.Range("L2").Value = 0.005536985397338867
$ws.Range("M2").Value = 0.01330232620239258
$ws.Range("B3").Value = 0.05794916152954101
$ws.Range("C3").Value = 0.02321267127990723
$ws.Range("D3").Value = 0.01250844001770019
$ws.Range("E3").Value = 0.01058416366577148
$ws.Range("F3").Value = 0.009735202789306641
$ws.Range("G3").Value = 0.002799844741821289
$ws.Range("H3").Value = 0.06732068061828614
$ws.Range("I3").Value = 0.02409334182739258
$ws.Range("J3").Value = 0.04278016090393066
$ws.Range("K3").Value = 0.0110142707824707
$ws.Range("L3").Value = 0.01713852882385254
$ws.Range("M3").Value = 0.01138029098510742
$ws.Range("B4").Value = 0.02223739624023437
$ws.Range("C4").Value = 0.01289124488830566
$ws.Range("D4").Value = 0.01388835906982422
$ws.Range("E4").Value = 0.001472616195678711
$ws.Range("F4").Value = 0.0301846981048584
$ws.Range("G4").Value = 0.007805442810058594
$ws.Range("H4").Value = 0.01461348533630371
$ws.Range("I4").Value = 0.01173267364501953
$ws.Range("J4").Value = 0.01222672462463379
$ws.Range("K4").Value = 0.008199977874755859
$ws.Range("L4").Value = 0.0352264404296875
$ws.Range("M4").Value = 0.008774614334106446
$ws.Range("B5").Value = 0.01409602165222168
$ws.Range("C5").Value = 0.01233382225036621
$ws.Range("D5").Value = 0.003233194351196289
$ws.Range("E5").Value = 0.01342692375183105
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.01256961822509766
$ws.Range("I5").Value = 0.01077876091003418
$ws.Range("J5").Value = 0.01663980484008789
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("B6").Value = 0.4370864391326904
$ws.Range("C6").Value = 0.01960210800170898
$ws.Range("D6").Value = 0.3234162330627441
$ws.Range("E6").Value = 0.01714372634887695
$ws.Range("F6").Value = 0.3555170059204101
$ws.Range("G6").Value = 0.01420173645019531
$ws.Range("H6").Value = 0.1269444942474365
$ws.Range("I6").Value = 0.01027321815490723
$ws.Range("J6").Value = 0.1189289569854736
$ws.Range("K6").Value = 0.01568608283996582
$ws.Range("L6").Value = 0.3811039924621582
$ws.Range("M6").Value = 0.01101088523864746
$ws.Range("B7").Value = 0.4893858909606933
$ws.Range("C7").Value = 0.05815262794494629
$ws.Range("D7").Value = 0.2168639659881592
$ws.Range("E7").Value = 0.01714091300964355
$ws.Range("F7").Value = 0.4836967945098877
$ws.Range("G7").Value = 0.01922645568847656
$ws.Range("H7").Value = 0.5476384162902832
$ws.Range("I7").Value = 0.0615788459777832
$ws.Range("J7").Value = 0.3476318836212158
$ws.Range("K7").Value = 0.03600363731384278
$ws.Range("L7").Value = 0.5488099098205567
$ws.Range("M7").Value = 0.01884818077087402